$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 40159
$ws.Range("I11").Value = 40159
$ws.Range("K11").Value = 40159
$ws.Range("M11").Value = -40019
$ws.Range("H18").Value = 394.56522
$ws.Range("I18").Value = 394.56522
$ws.Range("K18").Value = 394.56522
$ws.Range("M18").Value = -110.56522
$ws.Range("H33").Value = 8761.833000000001
$ws.Range("I33").Value = 14741.429
$ws.Range("J33").Value = 390.4
$ws.Range("K33").Value = 14741.429
$ws.Range("L33").Value = 390.4
$ws.Range("M33").Value = -14512.429
$ws.Range("N33").Value = -848.4
$ws.Range("H51").Value = 9099.200000000001
$ws.Range("I51").Value = 8570.571
$ws.Range("J51").Value = 10332.667
$ws.Range("K51").Value = 8570.571
$ws.Range("L51").Value = 10332.667
$ws.Range("M51").Value = -8086.571
$ws.Range("N51").Value = -11300.667
$ws.Range("H76").Value = 5433
$ws.Range("I76").Value = 6299
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 6299
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -5984
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 5433
$ws.Range("I79").Value = 6299
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 6299
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -5207
$ws.Range("N79").Value = -7184
$ws.Range("H116").Value = 8923.111000000001
$ws.Range("I116").Value = 4930.857
$ws.Range("J116").Value = 11463.637
$ws.Range("K116").Value = 4930.857
$ws.Range("L116").Value = 11463.637
$ws.Range("M116").Value = -1488.857
$ws.Range("N116").Value = -18347.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3169.463
$ws.Range("I2").Value = 2418.3953
$ws.Range("K2").Value = 2418.3953
$ws.Range("M2").Value = -2305.3953
$ws.Range("H61").Value = 5327.75
$ws.Range("J61").Value = 15332.667
$ws.Range("L61").Value = 15332.667
$ws.Range("N61").Value = -15756.667
$ws.Range("H62").Value = 50244
$ws.Range("J62").Value = 50244
$ws.Range("L62").Value = 50244
$ws.Range("N62").Value = -51492
$ws.Range("H65").Value = 50244
$ws.Range("J65").Value = 50244
$ws.Range("L65").Value = 150732
$ws.Range("N65").Value = -156972
$ws.Range("H116").Value = 3169.463
$ws.Range("I116").Value = 2418.3953
$ws.Range("K116").Value = 2418.3953
$ws.Range("M116").Value = -124.3953000000001
$ws.Range("H136").Value = 5327.75
$ws.Range("J136").Value = 15332.667
$ws.Range("L136").Value = 45998.001
$ws.Range("N136").Value = -51098.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3169.463
$ws.Range("I3").Value = 2418.3953
$ws.Range("K3").Value = 2418.3953
$ws.Range("M3").Value = -2304.3953

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5314.9585
$ws.Range("I58").Value = 2198.0588
$ws.Range("K58").Value = 2198.0588
$ws.Range("M58").Value = -1995.0588
$ws.Range("H122").Value = 2041.2307
$ws.Range("I122").Value = 2079.125
$ws.Range("J122").Value = 1980.6
$ws.Range("K122").Value = 6237.375
$ws.Range("L122").Value = 5941.799999999999
$ws.Range("M122").Value = -3787.375
$ws.Range("N122").Value = -10841.8
$ws.Range("H134").Value = 10458.091
$ws.Range("I134").Value = 9278.111000000001
$ws.Range("K134").Value = 27834.333
$ws.Range("M134").Value = -25299.333
$ws.Range("H136").Value = 5314.9585
$ws.Range("I136").Value = 2198.0588
$ws.Range("K136").Value = 6594.176399999999
$ws.Range("M136").Value = -4044.176399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3929.875
$ws.Range("I3").Value = 3929.875
$ws.Range("K3").Value = 11789.625
$ws.Range("M3").Value = -11677.625
$ws.Range("H56").Value = 6065.8
$ws.Range("I56").Value = 6065.8
$ws.Range("K56").Value = 6065.8
$ws.Range("M56").Value = -5535.8
$ws.Range("H68").Value = 1874.7
$ws.Range("J68").Value = 1860.7778
$ws.Range("L68").Value = 5582.3334
$ws.Range("N68").Value = -7204.3334
$ws.Range("H71").Value = 1874.7
$ws.Range("J71").Value = 1860.7778
$ws.Range("L71").Value = 16747.0002
$ws.Range("N71").Value = -24859.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 15684
$ws.Range("J15").Value = 15684
$ws.Range("L15").Value = 15684
$ws.Range("N15").Value = -16260
$ws.Range("H81").Value = 15684
$ws.Range("J81").Value = 15684
$ws.Range("L81").Value = 15684
$ws.Range("N81").Value = -17680
$ws.Range("H84").Value = 15684
$ws.Range("J84").Value = 15684
$ws.Range("L84").Value = 47052
$ws.Range("N84").Value = -57036
$ws.Range("H122").Value = 1415.625
$ws.Range("I122").Value = 1475
$ws.Range("K122").Value = 4425
$ws.Range("M122").Value = -1975

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2221.75
$ws.Range("I16").Value = 1696.3077
$ws.Range("J16").Value = 4498.6665
$ws.Range("K16").Value = 1696.3077
$ws.Range("L16").Value = 4498.6665
$ws.Range("M16").Value = -1526.3077
$ws.Range("N16").Value = -4838.6665
$ws.Range("H46").Value = 1823.3667
$ws.Range("I46").Value = 1853
$ws.Range("J46").Value = 1800.7059
$ws.Range("K46").Value = 1853
$ws.Range("L46").Value = 1800.7059
$ws.Range("M46").Value = -1665
$ws.Range("N46").Value = -2176.7059
$ws.Range("H69").Value = 3399999.2
$ws.Range("J69").Value = 3399999.2
$ws.Range("L69").Value = 3399999.2
$ws.Range("N69").Value = -3401621.2
$ws.Range("H72").Value = 3399999.2
$ws.Range("J72").Value = 3399999.2
$ws.Range("L72").Value = 10199997.6
$ws.Range("N72").Value = -10208109.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41622
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -128112
$ws.Range("H86").Value = 127499.5
$ws.Range("J86").Value = 127499.5
$ws.Range("L86").Value = 127499.5
$ws.Range("N86").Value = -129745.5
$ws.Range("H89").Value = 127499.5
$ws.Range("J89").Value = 127499.5
$ws.Range("L89").Value = 637497.5
$ws.Range("N89").Value = -648729.5
$ws.Range("H96").Value = 2876.7083
$ws.Range("I96").Value = 3544.4
$ws.Range("K96").Value = 3544.4
$ws.Range("M96").Value = -2171.4
$ws.Range("H107").Value = 1745.3043
$ws.Range("J107").Value = 1018
$ws.Range("L107").Value = 3054
$ws.Range("N107").Value = -6894
$ws.Range("H132").Value = 14163.5
$ws.Range("I132").Value = 11952.098
$ws.Range("J132").Value = 22406
$ws.Range("K132").Value = 35856.294
$ws.Range("L132").Value = 67218
$ws.Range("M132").Value = -33326.294
$ws.Range("N132").Value = -72278
